$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Fgf13 -> Scn5a -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf13"
$ws.Range("C2").Value = "Scn5a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.14444
$ws.Range("H2").Value = 0.43332
$ws.Range("I2").Value = 0.06801140868936309
$ws.Range("J2").Value = 0.06801140868936309
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06912333333333333
$ws.Range("N2").Value = 0.20737
$ws.Range("O2").Value = 0.01464591175868182
$ws.Range("P2").Value = 0.01464591175868182
$ws.Range("Q2").Value = 0.009984174266666666
$ws.Range("R2").Value = 0.0898575684
$ws.Range("S2").Value = 0.0009960890902480576
$ws.Range("T2").Value = 0.0009960890902480579

# Row 3: ECs -> Fgf13 -> Scn5a -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf13"
$ws.Range("C3").Value = "Scn5a"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.14444
$ws.Range("H3").Value = 0.43332
$ws.Range("I3").Value = 0.06801140868936309
$ws.Range("J3").Value = 0.06801140868936309
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.65051
$ws.Range("N3").Value = 13.95153
$ws.Range("O3").Value = 0.9853540882413181
$ws.Range("P3").Value = 0.9853540882413182
$ws.Range("Q3").Value = 0.6717196644
$ws.Range("R3").Value = 6.045476979600001
$ws.Range("S3").Value = 0.06701531959911503
$ws.Range("T3").Value = 0.06701531959911504

# Row 4: sCs -> Fgf13 -> Scn5a -> ECs
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Fgf13"
$ws.Range("C4").Value = "Scn5a"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.979321333333333
$ws.Range("H4").Value = 5.937964
$ws.Range("I4").Value = 0.9319885913106368
$ws.Range("J4").Value = 0.9319885913106369
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.06912333333333333
$ws.Range("N4").Value = 0.20737
$ws.Range("O4").Value = 0.01464591175868182
$ws.Range("P4").Value = 0.01464591175868182
$ws.Range("Q4").Value = 0.1368172882977778
$ws.Range("R4").Value = 1.23135559468
$ws.Range("S4").Value = 0.01364982266843376
$ws.Range("T4").Value = 0.01364982266843376

# Row 5: sCs -> Fgf13 -> Scn5a -> sCs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf13"
$ws.Range("C5").Value = "Scn5a"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.979321333333333
$ws.Range("H5").Value = 5.937964
$ws.Range("I5").Value = 0.9319885913106368
$ws.Range("J5").Value = 0.9319885913106369
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.65051
$ws.Range("N5").Value = 13.95153
$ws.Range("O5").Value = 0.9853540882413181
$ws.Range("P5").Value = 0.9853540882413182
$ws.Range("Q5").Value = 9.204853653879999
$ws.Range("R5").Value = 82.84368288492
$ws.Range("S5").Value = 0.918338768642203
$ws.Range("T5").Value = 0.9183387686422032
